$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 3036.625
$ws.Cells.Item(41, 9).Value = 3327.5715
$ws.Cells.Item(41, 10).Value = 1000
$ws.Cells.Item(41, 11).Value = 3327.5715
$ws.Cells.Item(41, 12).Value = 1000
$ws.Cells.Item(41, 13).Value = -2887.5715
$ws.Cells.Item(41, 14).Value = -1880

$ws.Cells.Item(43, 8).Value = 12996.5
$ws.Cells.Item(43, 10).Value = 19002
$ws.Cells.Item(43, 12).Value = 19002
$ws.Cells.Item(43, 14).Value = -19140

$ws.Cells.Item(57, 8).Value = 51722.5
$ws.Cells.Item(57, 10).Value = 51722.5
$ws.Cells.Item(57, 12).Value = 155167.5
$ws.Cells.Item(57, 14).Value = -156165.5

$ws.Cells.Item(76, 8).Value = 13387.333
$ws.Cells.Item(76, 9).Value = 13426
$ws.Cells.Item(76, 10).Value = 13252
$ws.Cells.Item(76, 11).Value = 13426
$ws.Cells.Item(76, 12).Value = 13252
$ws.Cells.Item(76, 13).Value = -13111
$ws.Cells.Item(76, 14).Value = -13882

$ws.Cells.Item(79, 8).Value = 13387.333
$ws.Cells.Item(79, 9).Value = 13426
$ws.Cells.Item(79, 10).Value = 13252
$ws.Cells.Item(79, 11).Value = 13426
$ws.Cells.Item(79, 12).Value = 13252
$ws.Cells.Item(79, 13).Value = -12334
$ws.Cells.Item(79, 14).Value = -15436

$ws.Cells.Item(113, 8).Value = 13781.6
$ws.Cells.Item(113, 10).Value = 13781.6
$ws.Cells.Item(113, 12).Value = 13781.6
$ws.Cells.Item(113, 14).Value = -20289.6

$ws.Cells.Item(116, 8).Value = 8249.786
$ws.Cells.Item(116, 9).Value = 9099.200000000001
$ws.Cells.Item(116, 11).Value = 9099.200000000001
$ws.Cells.Item(116, 13).Value = -5657.200000000001

$ws.Cells.Item(132, 8).Value = 1220.238
$ws.Cells.Item(132, 9).Value = 1220.238
$ws.Cells.Item(132, 11).Value = 3660.714
$ws.Cells.Item(132, 13).Value = -1130.714

$ws.Cells.Item(138, 8).Value = 3815.9707
$ws.Cells.Item(138, 10).Value = 3714.625
$ws.Cells.Item(138, 12).Value = 11143.875
$ws.Cells.Item(138, 14).Value = -21423.875

$ws.Cells.Item(141, 8).Value = 5720.9287
$ws.Cells.Item(141, 9).Value = 4021.4443
$ws.Cells.Item(141, 11).Value = 12064.3329
$ws.Cells.Item(141, 13).Value = -6884.332900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 14945.134
$ws.Cells.Item(2, 9).Value = 988.5714
$ws.Cells.Item(2, 10).Value = 27157.125
$ws.Cells.Item(2, 11).Value = 988.5714
$ws.Cells.Item(2, 12).Value = 27157.125
$ws.Cells.Item(2, 13).Value = -875.5714
$ws.Cells.Item(2, 14).Value = -27383.125

$ws.Cells.Item(5, 8).Value = 351.8889
$ws.Cells.Item(5, 9).Value = 168.5
$ws.Cells.Item(5, 11).Value = 168.5
$ws.Cells.Item(5, 13).Value = -56.5

$ws.Cells.Item(24, 8).Value = 27500
$ws.Cells.Item(24, 10).Value = 27500
$ws.Cells.Item(24, 12).Value = 27500
$ws.Cells.Item(24, 14).Value = -28248

$ws.Cells.Item(32, 8).Value = 2498.5557
$ws.Cells.Item(32, 9).Value = 1504.1621
$ws.Cells.Item(32, 11).Value = 1504.1621
$ws.Cells.Item(32, 13).Value = -1217.1621

$ws.Cells.Item(100, 8).Value = 27500
$ws.Cells.Item(100, 10).Value = 27500
$ws.Cells.Item(100, 12).Value = 27500
$ws.Cells.Item(100, 14).Value = -29664

$ws.Cells.Item(116, 8).Value = 14945.134
$ws.Cells.Item(116, 9).Value = 988.5714
$ws.Cells.Item(116, 10).Value = 27157.125
$ws.Cells.Item(116, 11).Value = 988.5714
$ws.Cells.Item(116, 12).Value = 27157.125
$ws.Cells.Item(116, 13).Value = 1305.4286
$ws.Cells.Item(116, 14).Value = -31745.125

$ws.Cells.Item(132, 8).Value = 3410.4119
$ws.Cells.Item(132, 9).Value = 1826.8462
$ws.Cells.Item(132, 11).Value = 5480.5386
$ws.Cells.Item(132, 13).Value = -2950.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 14945.134
$ws.Cells.Item(3, 9).Value = 988.5714
$ws.Cells.Item(3, 10).Value = 27157.125
$ws.Cells.Item(3, 11).Value = 988.5714
$ws.Cells.Item(3, 12).Value = 27157.125
$ws.Cells.Item(3, 13).Value = -874.5714
$ws.Cells.Item(3, 14).Value = -27385.125

$ws.Cells.Item(4, 8).Value = 351.8889
$ws.Cells.Item(4, 9).Value = 168.5
$ws.Cells.Item(4, 11).Value = 168.5
$ws.Cells.Item(4, 13).Value = -53.5

$ws.Cells.Item(86, 8).Value = 1965.1765
$ws.Cells.Item(86, 9).Value = 2185.2307
$ws.Cells.Item(86, 10).Value = 1250
$ws.Cells.Item(86, 11).Value = 2185.2307
$ws.Cells.Item(86, 12).Value = 1250
$ws.Cells.Item(86, 13).Value = -1062.2307
$ws.Cells.Item(86, 14).Value = -3496

$ws.Cells.Item(89, 8).Value = 1965.1765
$ws.Cells.Item(89, 9).Value = 2185.2307
$ws.Cells.Item(89, 10).Value = 1250
$ws.Cells.Item(89, 11).Value = 10926.1535
$ws.Cells.Item(89, 12).Value = 6250
$ws.Cells.Item(89, 13).Value = -5310.1535
$ws.Cells.Item(89, 14).Value = -17482

$ws.Cells.Item(94, 8).Value = 764.84
$ws.Cells.Item(94, 9).Value = 846.7895
$ws.Cells.Item(94, 10).Value = 505.33334
$ws.Cells.Item(94, 11).Value = 846.7895
$ws.Cells.Item(94, 12).Value = 505.33334
$ws.Cells.Item(94, 13).Value = -395.7895
$ws.Cells.Item(94, 14).Value = -1407.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 36743.605
$ws.Cells.Item(31, 9).Value = 3738.1177
$ws.Cells.Item(31, 11).Value = 3738.1177
$ws.Cells.Item(31, 13).Value = -3443.1177

$ws.Cells.Item(34, 8).Value = 36743.605
$ws.Cells.Item(34, 9).Value = 3738.1177
$ws.Cells.Item(34, 11).Value = 3738.1177
$ws.Cells.Item(34, 13).Value = -3536.1177

$ws.Cells.Item(99, 8).Value = 2499.75
$ws.Cells.Item(99, 9).Value = 2499.75
$ws.Cells.Item(99, 11).Value = 2499.75
$ws.Cells.Item(99, 13).Value = -1001.75

$ws.Cells.Item(105, 8).Value = 5072.5713
$ws.Cells.Item(105, 9).Value = 2899.4
$ws.Cells.Item(105, 11).Value = 2899.4
$ws.Cells.Item(105, 13).Value = -1152.4

$ws.Cells.Item(126, 8).Value = 2499.75
$ws.Cells.Item(126, 9).Value = 2499.75
$ws.Cells.Item(126, 11).Value = 7499.25
$ws.Cells.Item(126, 13).Value = -5029.25

$ws.Cells.Item(132, 8).Value = 4751
$ws.Cells.Item(132, 9).Value = 4574.2354
$ws.Cells.Item(132, 11).Value = 13722.7062
$ws.Cells.Item(132, 13).Value = -11192.7062

$ws.Cells.Item(134, 8).Value = 2538.4773
$ws.Cells.Item(134, 9).Value = 2053.647
$ws.Cells.Item(134, 10).Value = 4186.9
$ws.Cells.Item(134, 11).Value = 6160.941
$ws.Cells.Item(134, 12).Value = 12560.7
$ws.Cells.Item(134, 13).Value = -3625.941
$ws.Cells.Item(134, 14).Value = -17630.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 497
$ws.Cells.Item(97, 9).Value = 179.54546
$ws.Cells.Item(97, 11).Value = 538.6363799999999
$ws.Cells.Item(97, 13).Value = -42.63637999999992

$ws.Cells.Item(107, 8).Value = 711.7143
$ws.Cells.Item(107, 9).Value = 640.25
$ws.Cells.Item(107, 10).Value = 807
$ws.Cells.Item(107, 11).Value = 1920.75
$ws.Cells.Item(107, 12).Value = 2421
$ws.Cells.Item(107, 13).Value = -0.75
$ws.Cells.Item(107, 14).Value = -6261

$ws.Cells.Item(118, 8).Value = 2631.25
$ws.Cells.Item(118, 9).Value = 310
$ws.Cells.Item(118, 11).Value = 930
$ws.Cells.Item(118, 13).Value = 313

$ws.Cells.Item(132, 8).Value = 5675.4165
$ws.Cells.Item(132, 10).Value = 6215
$ws.Cells.Item(132, 12).Value = 55935
$ws.Cells.Item(132, 14).Value = -60995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 17992.125
$ws.Cells.Item(70, 9).Value = 5430.778
$ws.Cells.Item(70, 10).Value = 34142.43
$ws.Cells.Item(70, 11).Value = 5430.778
$ws.Cells.Item(70, 12).Value = 34142.43
$ws.Cells.Item(70, 13).Value = -5160.778
$ws.Cells.Item(70, 14).Value = -34682.43

$ws.Cells.Item(73, 8).Value = 17992.125
$ws.Cells.Item(73, 9).Value = 5430.778
$ws.Cells.Item(73, 10).Value = 34142.43
$ws.Cells.Item(73, 11).Value = 5430.778
$ws.Cells.Item(73, 12).Value = 34142.43
$ws.Cells.Item(73, 13).Value = -4494.778
$ws.Cells.Item(73, 14).Value = -36014.43

$ws.Cells.Item(80, 8).Value = 6198.4546
$ws.Cells.Item(80, 9).Value = 1863.3334
$ws.Cells.Item(80, 10).Value = 7824.125
$ws.Cells.Item(80, 11).Value = 1863.3334
$ws.Cells.Item(80, 12).Value = 7824.125
$ws.Cells.Item(80, 13).Value = -865.3334
$ws.Cells.Item(80, 14).Value = -9820.125

$ws.Cells.Item(83, 8).Value = 6198.4546
$ws.Cells.Item(83, 9).Value = 1863.3334
$ws.Cells.Item(83, 10).Value = 7824.125
$ws.Cells.Item(83, 11).Value = 9316.666999999999
$ws.Cells.Item(83, 12).Value = 39120.625
$ws.Cells.Item(83, 13).Value = -4324.666999999999
$ws.Cells.Item(83, 14).Value = -49104.625

$ws.Cells.Item(126, 8).Value = 2738.6843
$ws.Cells.Item(126, 9).Value = 1508.1333
$ws.Cells.Item(126, 10).Value = 7353.25
$ws.Cells.Item(126, 11).Value = 4524.3999
$ws.Cells.Item(126, 12).Value = 22059.75
$ws.Cells.Item(126, 13).Value = -2054.3999
$ws.Cells.Item(126, 14).Value = -26999.75

$ws.Cells.Item(132, 8).Value = 132644.62
$ws.Cells.Item(132, 9).Value = 251785.75
$ws.Cells.Item(132, 10).Value = 13503.5
$ws.Cells.Item(132, 11).Value = 755357.25
$ws.Cells.Item(132, 12).Value = 40510.5
$ws.Cells.Item(132, 13).Value = -752827.25
$ws.Cells.Item(132, 14).Value = -45570.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4880.0557
$ws.Cells.Item(7, 9).Value = 2582.1785
$ws.Cells.Item(7, 11).Value = 2582.1785
$ws.Cells.Item(7, 13).Value = -2470.1785

$ws.Cells.Item(64, 8).Value = 19333.334
$ws.Cells.Item(64, 10).Value = 24000
$ws.Cells.Item(64, 12).Value = 24000
$ws.Cells.Item(64, 14).Value = -24450

$ws.Cells.Item(67, 8).Value = 19333.334
$ws.Cells.Item(67, 10).Value = 24000
$ws.Cells.Item(67, 12).Value = 24000
$ws.Cells.Item(67, 14).Value = -25560

$ws.Cells.Item(126, 8).Value = 4880.0557
$ws.Cells.Item(126, 9).Value = 2582.1785
$ws.Cells.Item(126, 11).Value = 7746.5355
$ws.Cells.Item(126, 13).Value = -5276.5355

$ws.Cells.Item(132, 8).Value = 4410.2173
$ws.Cells.Item(132, 9).Value = 2853
$ws.Cells.Item(132, 10).Value = 6832.5557
$ws.Cells.Item(132, 11).Value = 8559
$ws.Cells.Item(132, 12).Value = 20497.6671
$ws.Cells.Item(132, 13).Value = -6029
$ws.Cells.Item(132, 14).Value = -25557.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 80000
$ws.Cells.Item(63, 10).Value = 80000
$ws.Cells.Item(63, 12).Value = 80000
$ws.Cells.Item(63, 14).Value = -81248

$ws.Cells.Item(66, 8).Value = 80000
$ws.Cells.Item(66, 10).Value = 80000
$ws.Cells.Item(66, 12).Value = 240000
$ws.Cells.Item(66, 14).Value = -246240

$ws.Cells.Item(126, 8).Value = 1575.1578
$ws.Cells.Item(126, 9).Value = 960.6667
$ws.Cells.Item(126, 10).Value = 3879.5
$ws.Cells.Item(126, 11).Value = 2882.0001
$ws.Cells.Item(126, 12).Value = 11638.5
$ws.Cells.Item(126, 13).Value = -412.0001000000002
$ws.Cells.Item(126, 14).Value = -16578.5

$ws.Cells.Item(132, 8).Value = 6066.4287
$ws.Cells.Item(132, 9).Value = 5846.5293
$ws.Cells.Item(132, 10).Value = 7001
$ws.Cells.Item(132, 11).Value = 17539.5879
$ws.Cells.Item(132, 12).Value = 21003
$ws.Cells.Item(132, 13).Value = -15009.5879
$ws.Cells.Item(132, 14).Value = -26063
